# Updated symbol list on Sat Jan  7 05:52:34 UTC 2023 with GitHub Actions
#
# This refreshes the "Price" (column D) and "Volume(1h)" (column E) figures
# for the crypto-ranking rows on Sheet1, matching a new scrape snapshot.
# Every value is textual (e.g. "261.17", "1.71%") exactly like the source
# data, so each cell is forced to Text before the write (NumberFormat "@")
# and then the style is put back to "Normal" so no extra formatting /
# style index lingers on the cell - only the literal text changes, just
# like the upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# row 2 - BNB
Set-TextValue "D2" "261.17"
Set-TextValue "E2" "1.71%"

# row 3 - OKB
Set-TextValue "D3" "27.48"
Set-TextValue "E3" "1.68%"

# row 4 - HuobiToken
Set-TextValue "D4" "4.767"
Set-TextValue "E4" "10.33%"

# row 5 - Cronos
Set-TextValue "D5" "0.06073"
Set-TextValue "E5" "3.15%"

# row 6 - KuCoinToken
Set-TextValue "D6" "6.661"
Set-TextValue "E6" "0.41%"

# row 7 - MXToken
Set-TextValue "D7" "0.8606"
Set-TextValue "E7" "1.21%"

# row 8 - FTXToken
Set-TextValue "D8" "0.9246"
Set-TextValue "E8" "-1.35%"

# row 9 - WazirX
Set-TextValue "D9" "0.1406"
Set-TextValue "E9" "1.53%"

# row 10 - LiechtensteinCryptoassetsExchange
Set-TextValue "D10" "0.04923"
Set-TextValue "E10" "3.71%"

# row 11 - MandalaExchangeToken
Set-TextValue "D11" "0.07100"
Set-TextValue "E11" "0.26%"

# row 12 - BitrueCoin
Set-TextValue "D12" "0.03053"
Set-TextValue "E12" "-0.73%"

# row 13 - BitMartToken
Set-TextValue "D13" "0.09083"
Set-TextValue "E13" "-0.37%"

# row 14 - BitForexToken
Set-TextValue "D14" "0.001531"
Set-TextValue "E14" "0.32%"

# row 15 - One
Set-TextValue "D15" "0.0006093"
Set-TextValue "E15" "0.68%"

# row 16 - TigerCash
Set-TextValue "D16" "0.006078"
Set-TextValue "E16" "-1.03%"

# row 17 - LEO
Set-TextValue "D17" "3.454"
Set-TextValue "E17" "-1.10%"

# row 18 - GateToken
Set-TextValue "D18" "3.160"
Set-TextValue "E18" "-0.71%"

# row 19 - BTSEToken (price unchanged)
Set-TextValue "E19" "-2.23%"

# row 21 - ProBitToken
Set-TextValue "D21" "0.1297"
Set-TextValue "E21" "2.12%"

# row 22 - MCDex
Set-TextValue "D22" "4.113"
Set-TextValue "E22" "4.89%"

# row 23 - CoinExToken
Set-TextValue "D23" "0.04262"
Set-TextValue "E23" "0.23%"

# row 24 - BitKan (price unchanged)
Set-TextValue "E24" "0.02%"

# row 25 - HotbitToken (price unchanged)
Set-TextValue "E25" "-8.66%"

# row 26 - NitroEx (volume unchanged)
Set-TextValue "D26" "0.0001200"

# row 27 - UpBots (price unchanged)
Set-TextValue "E27" "3.10%"

# row 40 - IDEX
Set-TextValue "D40" "0.03879"
Set-TextValue "E40" "1.63%"

# row 41 - BKEXToken (price unchanged)
Set-TextValue "E41" "1.46%"

# row 42 - KickToken
Set-TextValue "D42" "0.004131"
Set-TextValue "E42" "-33.49%"

# row 43 - LocalTraders
Set-TextValue "D43" "0.01501"
Set-TextValue "E43" "6.88%"

# row 44 - CEJI (price unchanged)
Set-TextValue "E44" "0.56%"

# row 45 - CoinLion
Set-TextValue "D45" "0.00005299"
Set-TextValue "E45" "-1.26%"

# row 46 - Kangarootoken (price unchanged)
Set-TextValue "E46" "0.02%"

# row 47 - CoinbaseStockToken (price unchanged)
Set-TextValue "E47" "-17.21%"

# row 48 - BOLO
Set-TextValue "D48" "0.1353"
Set-TextValue "E48" "-46.38%"

# row 49 - CryptobidCoin (price unchanged)
Set-TextValue "E49" "0.02%"

# row 50 - SpecialPowerGold (price unchanged)
Set-TextValue "E50" "0.02%"

Write-Output "Applied 61 price/volume cell updates"
